$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 & 6: update the latest/next period dates and ILR source links.
# Order matches the shared-string append order of the target file:
# new period text first, then the two distinct ILR links.
$ws.Range("D5").Value = "Aug 2024 – Jul 2025 (Nov 25)"
$ws.Range("D6").Value = "Aug 2024 – Jul 2025 (Nov 25)"

$ws.Range("C5").Value = "Aug 2023 – Jul 2024 (28/11/24)"
$ws.Range("C6").Value = "Aug 2023 – Jul 2024 (28/11/24)"

$ws.Range("B5").Value = "<a href='https://explore-education-statistics.service.gov.uk/data-catalogue/data-set/b930498d-b4f0-416d-a086-7acee1be8179'>Individualised Learner Record</a>"
$ws.Range("B6").Value = "<a href='https://explore-education-statistics.service.gov.uk/data-tables/permalink/47e8710e-ccb7-401c-ab25-08dd34489990'>Individualised Learner Record</a>"

# Reset style on C5/C6 back to default (General) -- the s="1" attribute is dropped in the edit
$ws.Range("C5").Style = "Normal"
$ws.Range("C6").Style = "Normal"

# Update the active selection to B7
$ws.Range("B7").Select()
